$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = "61.247.62"
$ws.Cells.Item(2, 5).Value = "  -4.45%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.307.31"
$ws.Cells.Item(3, 5).Value = "  -5.01%  "

# Row 5
Set-TextValue $ws 5 4 "566.51"
$ws.Cells.Item(5, 5).Value = "  -3.30%  "

# Row 6
Set-TextValue $ws 6 4 "128.14"
$ws.Cells.Item(6, 5).Value = "  -2.99%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.308.32"
$ws.Cells.Item(8, 5).Value = "  -4.96%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -1.21%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -4.30%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -4.13%  "

# Row 12
Set-TextValue $ws 12 4 "0.375"
$ws.Cells.Item(12, 5).Value = "  -3.11%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "3.872.67"
$ws.Cells.Item(13, 5).Value = "  -5.04%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -0.60%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.305.60"
$ws.Cells.Item(15, 5).Value = "  -5.13%  "

# Row 16
Set-TextValue $ws 16 4 "0.0000167"
$ws.Cells.Item(16, 5).Value = "  -5.68%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "61.341.61"
$ws.Cells.Item(17, 5).Value = "  -4.35%  "

# Row 18
Set-TextValue $ws 18 4 "24.26"
$ws.Cells.Item(18, 5).Value = "  -0.45%  "

# Row 19
Set-TextValue $ws 19 4 "5.66"
$ws.Cells.Item(19, 5).Value = "  -1.14%  "

# Row 20
Set-TextValue $ws 20 4 "13.39"
$ws.Cells.Item(20, 5).Value = "  -0.84%  "

# Row 21
Set-TextValue $ws 21 4 "8.94"
$ws.Cells.Item(21, 5).Value = "  -10.50%  "

# Row 22
Set-TextValue $ws 22 4 "354.80"
$ws.Cells.Item(22, 5).Value = "  -7.71%  "

# Row 23
Set-TextValue $ws 23 4 "0.553"
$ws.Cells.Item(23, 5).Value = "  -3.81%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -0.09%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "3.437.68"
$ws.Cells.Item(25, 5).Value = "  -5.10%  "

# Row 26
Set-TextValue $ws 26 4 "69.00"
$ws.Cells.Item(26, 5).Value = "  -7.69%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -5.02%  "

# Row 28
Set-TextValue $ws 28 4 "0.999"
$ws.Cells.Item(28, 5).Value = "  -0.20%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.30%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.16%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws 31 4 "2.11"
$ws.Cells.Item(31, 5).Value = "  -5.38%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws 32 4 "7.79"
$ws.Cells.Item(32, 5).Value = "  -1.75%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.01%  "

# Row 34
Set-TextValue $ws 34 4 "0.148"
$ws.Cells.Item(34, 5).Value = "  -3.05%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "3.335.04"
$ws.Cells.Item(35, 5).Value = "  -5.04%  "

# Row 36
Set-TextValue $ws 36 4 "22.62"
$ws.Cells.Item(36, 5).Value = "  -1.48%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +1.95%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -0.31%  "

# Row 39
Set-TextValue $ws 39 4 "162.61"
$ws.Cells.Item(39, 5).Value = "  -0.40%  "

# Row 40
Set-TextValue $ws 40 4 "1.47"
$ws.Cells.Item(40, 5).Value = "  -3.23%  "

# Row 41
Set-TextValue $ws 41 4 "0.0754"
$ws.Cells.Item(41, 5).Value = "  -3.27%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -0.10%  "

# Row 43
Set-TextValue $ws 43 4 "4.39"
$ws.Cells.Item(43, 5).Value = "  +1.47%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.99%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -7.28%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -1.67%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -4.31%  "

# Row 48
Set-TextValue $ws 48 4 "22.16"
$ws.Cells.Item(48, 5).Value = "  -8.21%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.96%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -8.02%  "

# Row 51
Set-TextValue $ws 51 4 "21.24"
$ws.Cells.Item(51, 5).Value = "  +3.29%  "
